$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 10 new rows before the existing row 972, pushing the old
# rows 972:985 down to 982:995 (same behaviour as in the target workbook).
$ws.Rows("972:981").Insert()

# Fill in the 10 newly inserted rows (972:981) with the new weekly data.
# Columns A,B,C,E,F,G,N,Q,R keep the constant values used throughout the sheet.

# Row 972
$ws.Range("A972").Value = 6
$ws.Range("B972").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C972").Value = "Metropolitana"
$ws.Range("D972").Value = 44448
$ws.Range("E972").Value = 13
$ws.Range("F972").Value = 100112045
$ws.Range("G972").Value = "Zapallo"
$ws.Range("H972").Value = "Camote"
$ws.Range("I972").Value = "1a (guarda)"
$ws.Range("J972").Value = 900
$ws.Range("K972").Value = 500
$ws.Range("L972").Value = 550
$ws.Range("M972").Value = 533
$ws.Range("N972").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O972").Value = "Provincia de Maipo"
$ws.Range("P972").Value = 533
$ws.Range("Q972").Value = 1
$ws.Range("R972").Value = "Hortaliza"

# Row 973
$ws.Range("A973").Value = 6
$ws.Range("B973").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C973").Value = "Metropolitana"
$ws.Range("D973").Value = 44448
$ws.Range("E973").Value = 13
$ws.Range("F973").Value = 100112045
$ws.Range("G973").Value = "Zapallo"
$ws.Range("H973").Value = "Camote"
$ws.Range("I973").Value = "1a (guarda)"
$ws.Range("J973").Value = 1100
$ws.Range("K973").Value = 500
$ws.Range("L973").Value = 550
$ws.Range("M973").Value = 532
$ws.Range("N973").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O973").Value = "Provincia de Maipo"
$ws.Range("P973").Value = 532
$ws.Range("Q973").Value = 1
$ws.Range("R973").Value = "Hortaliza"

# Row 974
$ws.Range("A974").Value = 6
$ws.Range("B974").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C974").Value = "Metropolitana"
$ws.Range("D974").Value = 44448
$ws.Range("E974").Value = 13
$ws.Range("F974").Value = 100112045
$ws.Range("G974").Value = "Zapallo"
$ws.Range("H974").Value = "Camote"
$ws.Range("I974").Value = "1a nueva(o)"
$ws.Range("J974").Value = 1600
$ws.Range("K974").Value = 550
$ws.Range("L974").Value = 550
$ws.Range("M974").Value = 550
$ws.Range("N974").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O974").Value = "Perú"
$ws.Range("P974").Value = 550
$ws.Range("Q974").Value = 1
$ws.Range("R974").Value = "Hortaliza"

# Row 975
$ws.Range("A975").Value = 6
$ws.Range("B975").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C975").Value = "Metropolitana"
$ws.Range("D975").Value = 44448
$ws.Range("E975").Value = 13
$ws.Range("F975").Value = 100112045
$ws.Range("G975").Value = "Zapallo"
$ws.Range("H975").Value = "Camote"
$ws.Range("I975").Value = "2a (guarda)"
$ws.Range("J975").Value = 1600
$ws.Range("K975").Value = 250
$ws.Range("L975").Value = 300
$ws.Range("M975").Value = 278
$ws.Range("N975").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O975").Value = "Provincia de Maipo"
$ws.Range("P975").Value = 278
$ws.Range("Q975").Value = 1
$ws.Range("R975").Value = "Hortaliza"

# Row 976
$ws.Range("A976").Value = 6
$ws.Range("B976").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C976").Value = "Metropolitana"
$ws.Range("D976").Value = 44448
$ws.Range("E976").Value = 13
$ws.Range("F976").Value = 100112045
$ws.Range("G976").Value = "Zapallo"
$ws.Range("H976").Value = "Camote"
$ws.Range("I976").Value = "2a (guarda)"
$ws.Range("J976").Value = 2000
$ws.Range("K976").Value = 250
$ws.Range("L976").Value = 300
$ws.Range("M976").Value = 280
$ws.Range("N976").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O976").Value = "Provincia de Maipo"
$ws.Range("P976").Value = 280
$ws.Range("Q976").Value = 1
$ws.Range("R976").Value = "Hortaliza"

# Row 977
$ws.Range("A977").Value = 6
$ws.Range("B977").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C977").Value = "Metropolitana"
$ws.Range("D977").Value = 44448
$ws.Range("E977").Value = 13
$ws.Range("F977").Value = 100112045
$ws.Range("G977").Value = "Zapallo"
$ws.Range("H977").Value = "Camote"
$ws.Range("I977").Value = "2a nueva(o)"
$ws.Range("J977").Value = 1100
$ws.Range("K977").Value = 400
$ws.Range("L977").Value = 400
$ws.Range("M977").Value = 400
$ws.Range("N977").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O977").Value = "Perú"
$ws.Range("P977").Value = 400
$ws.Range("Q977").Value = 1
$ws.Range("R977").Value = "Hortaliza"

# Row 978
$ws.Range("A978").Value = 6
$ws.Range("B978").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C978").Value = "Metropolitana"
$ws.Range("D978").Value = 44448
$ws.Range("E978").Value = 13
$ws.Range("F978").Value = 100112045
$ws.Range("G978").Value = "Zapallo"
$ws.Range("H978").Value = "Camote"
$ws.Range("I978").Value = "3a (guarda)"
$ws.Range("J978").Value = 300
$ws.Range("K978").Value = 150
$ws.Range("L978").Value = 150
$ws.Range("M978").Value = 150
$ws.Range("N978").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O978").Value = "Provincia de Maipo"
$ws.Range("P978").Value = 150
$ws.Range("Q978").Value = 1
$ws.Range("R978").Value = "Hortaliza"

# Row 979
$ws.Range("A979").Value = 6
$ws.Range("B979").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C979").Value = "Metropolitana"
$ws.Range("D979").Value = 44448
$ws.Range("E979").Value = 13
$ws.Range("F979").Value = 100112045
$ws.Range("G979").Value = "Zapallo"
$ws.Range("H979").Value = "Camote"
$ws.Range("I979").Value = "3a (guarda)"
$ws.Range("J979").Value = 400
$ws.Range("K979").Value = 150
$ws.Range("L979").Value = 150
$ws.Range("M979").Value = 150
$ws.Range("N979").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O979").Value = "Región de O'Higgins"
$ws.Range("P979").Value = 150
$ws.Range("Q979").Value = 1
$ws.Range("R979").Value = "Hortaliza"

# Row 980
$ws.Range("A980").Value = 6
$ws.Range("B980").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C980").Value = "Metropolitana"
$ws.Range("D980").Value = 44448
$ws.Range("E980").Value = 13
$ws.Range("F980").Value = 100112045
$ws.Range("G980").Value = "Zapallo"
$ws.Range("H980").Value = "Paine"
$ws.Range("I980").Value = "1a (guarda)"
$ws.Range("J980").Value = 3500
$ws.Range("K980").Value = 170
$ws.Range("L980").Value = 170
$ws.Range("M980").Value = 170
$ws.Range("N980").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O980").Value = "Región de O'Higgins"
$ws.Range("P980").Value = 170
$ws.Range("Q980").Value = 1
$ws.Range("R980").Value = "Hortaliza"

# Row 981
$ws.Range("A981").Value = 6
$ws.Range("B981").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C981").Value = "Metropolitana"
$ws.Range("D981").Value = 44448
$ws.Range("E981").Value = 13
$ws.Range("F981").Value = 100112045
$ws.Range("G981").Value = "Zapallo"
$ws.Range("H981").Value = "Paine"
$ws.Range("I981").Value = "2a (guarda)"
$ws.Range("J981").Value = 2200
$ws.Range("K981").Value = 120
$ws.Range("L981").Value = 120
$ws.Range("M981").Value = 120
$ws.Range("N981").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O981").Value = "Región de O'Higgins"
$ws.Range("P981").Value = 120
$ws.Range("Q981").Value = 1
$ws.Range("R981").Value = "Hortaliza"
